$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 352.84
$ws.Range("I33").Value = 313.9565
$ws.Range("K33").Value = 313.9565
$ws.Range("M33").Value = -84.95650000000001
$ws.Range("H100").Value = 2780.8333
$ws.Range("I100").Value = 2258
$ws.Range("J100").Value = 3154.2856
$ws.Range("K100").Value = 2258
$ws.Range("L100").Value = 3154.2856
$ws.Range("M100").Value = -1717
$ws.Range("N100").Value = -4236.2856
$ws.Range("H108").Value = 30309.25
$ws.Range("J108").Value = 30309.25
$ws.Range("L108").Value = 30309.25
$ws.Range("N108").Value = -37989.25
$ws.Range("H110").Value = 45700
$ws.Range("J110").Value = 45700
$ws.Range("L110").Value = 45700
$ws.Range("N110").Value = -53880
$ws.Range("H116").Value = 4521.2856
$ws.Range("I116").Value = 2497.2
$ws.Range("K116").Value = 2497.2
$ws.Range("M116").Value = 944.8000000000002
$ws.Range("H125").Value = 1091.6364
$ws.Range("I125").Value = 285.6
$ws.Range("J125").Value = 1763.3334
$ws.Range("K125").Value = 2570.4
$ws.Range("L125").Value = 15870.0006
$ws.Range("M125").Value = -110.4000000000001
$ws.Range("N125").Value = -20790.0006
$ws.Range("H137").Value = 1932.9166
$ws.Range("I137").Value = 1815.7894
$ws.Range("K137").Value = 5447.3682
$ws.Range("M137").Value = -2897.3682

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2154.4707
$ws.Range("I2").Value = 1443.909
$ws.Range("J2").Value = 3457.1667
$ws.Range("K2").Value = 1443.909
$ws.Range("L2").Value = 3457.1667
$ws.Range("M2").Value = -1330.909
$ws.Range("N2").Value = -3683.1667
$ws.Range("H11").Value = 1003
$ws.Range("I11").Value = 1003
$ws.Range("K11").Value = 1003
$ws.Range("M11").Value = -859
$ws.Range("H45").Value = 3228.4827
$ws.Range("I45").Value = 3519.2727
$ws.Range("J45").Value = 3050.7778
$ws.Range("K45").Value = 3519.2727
$ws.Range("L45").Value = 3050.7778
$ws.Range("M45").Value = -3142.2727
$ws.Range("N45").Value = -3804.7778
$ws.Range("H116").Value = 2154.4707
$ws.Range("I116").Value = 1443.909
$ws.Range("J116").Value = 3457.1667
$ws.Range("K116").Value = 1443.909
$ws.Range("L116").Value = 3457.1667
$ws.Range("M116").Value = 850.0909999999999
$ws.Range("N116").Value = -8045.1667
$ws.Range("H122").Value = 1766.2
$ws.Range("I122").Value = 1766.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5298.6
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2848.6
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 40429
$ws.Range("J123").Value = 40429
$ws.Range("L123").Value = 40429
$ws.Range("N123").Value = -50229

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2154.4707
$ws.Range("I3").Value = 1443.909
$ws.Range("J3").Value = 3457.1667
$ws.Range("K3").Value = 1443.909
$ws.Range("L3").Value = 3457.1667
$ws.Range("M3").Value = -1329.909
$ws.Range("N3").Value = -3685.1667
$ws.Range("H12").Value = 5200
$ws.Range("I12").Value = 400
$ws.Range("J12").Value = 10000
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = -232
$ws.Range("N12").Value = -10336
$ws.Range("H99").Value = 1741.5834
$ws.Range("I99").Value = 1495.4
$ws.Range("K99").Value = 1495.4
$ws.Range("M99").Value = 2.599999999999909
$ws.Range("H107").Value = 1600.375
$ws.Range("I107").Value = 756
$ws.Range("J107").Value = 3007.6667
$ws.Range("K107").Value = 756
$ws.Range("L107").Value = 3007.6667
$ws.Range("M107").Value = 1164
$ws.Range("N107").Value = -6847.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4484.7827
$ws.Range("I99").Value = 3272.7646
$ws.Range("J99").Value = 7918.8335
$ws.Range("K99").Value = 3272.7646
$ws.Range("L99").Value = 7918.8335
$ws.Range("M99").Value = -1774.7646
$ws.Range("N99").Value = -10914.8335
$ws.Range("H122").Value = 1296
$ws.Range("I122").Value = 1237.25
$ws.Range("J122").Value = 1325.375
$ws.Range("K122").Value = 3711.75
$ws.Range("L122").Value = 3976.125
$ws.Range("M122").Value = -1261.75
$ws.Range("N122").Value = -8876.125
$ws.Range("H126").Value = 4484.7827
$ws.Range("I126").Value = 3272.7646
$ws.Range("J126").Value = 7918.8335
$ws.Range("K126").Value = 9818.293799999999
$ws.Range("L126").Value = 23756.5005
$ws.Range("M126").Value = -7348.293799999999
$ws.Range("N126").Value = -28696.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 708.5
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H60").Value = 300
$ws.Range("I60").Value = 300
$ws.Range("K60").Value = 900
$ws.Range("M60").Value = -649
$ws.Range("H131").Value = 110692.19
$ws.Range("I131").Value = 807.5
$ws.Range("J131").Value = 115744.36
$ws.Range("K131").Value = 2422.5
$ws.Range("L131").Value = 347233.08
$ws.Range("M131").Value = 2617.5
$ws.Range("N131").Value = -357313.08

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 4201221.5
$ws.Range("I14").Value = 4201221.5
$ws.Range("K14").Value = 4201221.5
$ws.Range("M14").Value = -4201053.5
$ws.Range("H122").Value = 1882
$ws.Range("I122").Value = 1920.2
$ws.Range("K122").Value = 5760.6
$ws.Range("M122").Value = -3310.6
$ws.Range("H126").Value = 4857.1665
$ws.Range("I126").Value = 4452.636
$ws.Range("J126").Value = 5492.857
$ws.Range("K126").Value = 13357.908
$ws.Range("L126").Value = 16478.571
$ws.Range("M126").Value = -10887.908
$ws.Range("N126").Value = -21418.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H122").Value = 1723.6
$ws.Range("I122").Value = 1514.25
$ws.Range("K122").Value = 4542.75
$ws.Range("M122").Value = -2092.75
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = 0
$ws.Range("H126").Value = 1290.4286
$ws.Range("I126").Value = 755.5
$ws.Range("K126").Value = 2266.5
$ws.Range("M126").Value = 203.5
